# "Added Jan 28 Tournament" - appends 14 new ELO rows (rows 62-75,
# tournament entries #61-#74) to the bottom of the existing time table on Sheet1,
# growing the used range from A1:F61 to A1:F75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the bold/bordered/centered "entry number" index (style used by
# A2:A61 already). Clone that formatting down onto the new rows before filling
# in values, so A62:A75 pick up the same style as the rows above them.
$ws.Range("A61").Copy()
$ws.Range("A62:A75").PasteSpecial(-4122)

# Row 62 - tournament entry #61
$ws.Cells.Item(62, 1).Value = 61
$ws.Range("B62").Value = "'"
$ws.Range("B62").Style = "Normal"
$ws.Range("C62").Value = 1227.640338358163
$ws.Range("D62").Value = "'"
$ws.Range("D62").Style = "Normal"
$ws.Range("E62").Value = "'"
$ws.Range("E62").Style = "Normal"
$ws.Range("F62").Value = 1205.101394671996

# Row 63 - tournament entry #62
$ws.Cells.Item(63, 1).Value = 62
$ws.Range("B63").Value = 1300.792966980156
$ws.Range("C63").Value = "'"
$ws.Range("C63").Style = "Normal"
$ws.Range("D63").Value = "'"
$ws.Range("D63").Style = "Normal"
$ws.Range("E63").Value = 1173.622844683383
$ws.Range("F63").Value = "'"
$ws.Range("F63").Style = "Normal"

# Row 64 - tournament entry #63
$ws.Cells.Item(64, 1).Value = 63
$ws.Range("B64").Value = 1306.939512229892
$ws.Range("C64").Value = 1210.231062826311
$ws.Range("D64").Value = "'"
$ws.Range("D64").Style = "Normal"
$ws.Range("E64").Value = "'"
$ws.Range("E64").Style = "Normal"
$ws.Range("F64").Value = "'"
$ws.Range("F64").Style = "Normal"

# Row 65 - tournament entry #64
$ws.Cells.Item(65, 1).Value = 64
$ws.Range("B65").Value = "'"
$ws.Range("B65").Style = "Normal"
$ws.Range("C65").Value = "'"
$ws.Range("C65").Style = "Normal"
$ws.Range("D65").Value = "'"
$ws.Range("D65").Style = "Normal"
$ws.Range("E65").Value = 1167.476299433647
$ws.Range("F65").Value = 1222.510670203849

# Row 66 - tournament entry #65
$ws.Cells.Item(66, 1).Value = 65
$ws.Range("B66").Value = "'"
$ws.Range("B66").Style = "Normal"
$ws.Range("C66").Value = 1233.268980622833
$ws.Range("D66").Value = "'"
$ws.Range("D66").Style = "Normal"
$ws.Range("E66").Value = 1156.524349367828
$ws.Range("F66").Value = "'"
$ws.Range("F66").Style = "Normal"

# Row 67 - tournament entry #66
$ws.Cells.Item(67, 1).Value = 66
$ws.Range("B67").Value = 1283.90159443337
$ws.Range("C67").Value = "'"
$ws.Range("C67").Style = "Normal"
$ws.Range("D67").Value = "'"
$ws.Range("D67").Style = "Normal"
$ws.Range("E67").Value = "'"
$ws.Range("E67").Style = "Normal"
$ws.Range("F67").Value = 1233.462620269667

# Row 68 - tournament entry #67
$ws.Cells.Item(68, 1).Value = 67
$ws.Range("B68").Value = "'"
$ws.Range("B68").Style = "Normal"
$ws.Range("C68").Value = 1242.596162400927
$ws.Range("D68").Value = "'"
$ws.Range("D68").Style = "Normal"
$ws.Range("E68").Value = "'"
$ws.Range("E68").Style = "Normal"
$ws.Range("F68").Value = 1222.150502193288

# Row 69 - tournament entry #68
$ws.Cells.Item(69, 1).Value = 68
$ws.Range("B69").Value = 1295.213712509749
$ws.Range("C69").Value = "'"
$ws.Range("C69").Style = "Normal"
$ws.Range("D69").Value = "'"
$ws.Range("D69").Style = "Normal"
$ws.Range("E69").Value = 1147.197167589735
$ws.Range("F69").Value = "'"
$ws.Range("F69").Style = "Normal"

# Row 70 - tournament entry #69
$ws.Cells.Item(70, 1).Value = 69
$ws.Range("B70").Value = 1300.308918504197
$ws.Range("C70").Value = 1225.362651103099
$ws.Range("D70").Value = "'"
$ws.Range("D70").Style = "Normal"
$ws.Range("E70").Value = "'"
$ws.Range("E70").Style = "Normal"
$ws.Range("F70").Value = "'"
$ws.Range("F70").Style = "Normal"

# Row 71 - tournament entry #70
$ws.Cells.Item(71, 1).Value = 70
$ws.Range("B71").Value = "'"
$ws.Range("B71").Style = "Normal"
$ws.Range("C71").Value = "'"
$ws.Range("C71").Style = "Normal"
$ws.Range("D71").Value = "'"
$ws.Range("D71").Style = "Normal"
$ws.Range("E71").Value = 1142.101961595287
$ws.Range("F71").Value = 1239.384013491115

# Row 72 - tournament entry #71
$ws.Cells.Item(72, 1).Value = 71
$ws.Range("B72").Value = "'"
$ws.Range("B72").Style = "Normal"
$ws.Range("C72").Value = 1215.906131994833
$ws.Range("D72").Value = "'"
$ws.Range("D72").Style = "Normal"
$ws.Range("E72").Value = 1134.177109754972
$ws.Range("F72").Value = "'"
$ws.Range("F72").Style = "Normal"

# Row 73 - tournament entry #72
$ws.Cells.Item(73, 1).Value = 72
$ws.Range("B73").Value = 1309.765437612463
$ws.Range("C73").Value = "'"
$ws.Range("C73").Style = "Normal"
$ws.Range("D73").Value = "'"
$ws.Range("D73").Style = "Normal"
$ws.Range("E73").Value = "'"
$ws.Range("E73").Style = "Normal"
$ws.Range("F73").Value = 1247.308865331431

# Row 74 - tournament entry #73
$ws.Cells.Item(74, 1).Value = 73
$ws.Range("B74").Value = 1289.146991311059
$ws.Range("C74").Value = "'"
$ws.Range("C74").Style = "Normal"
$ws.Range("D74").Value = "'"
$ws.Range("D74").Style = "Normal"
$ws.Range("E74").Value = "'"
$ws.Range("E74").Style = "Normal"
$ws.Range("F74").Value = 1267.927311632835

# Row 75 - tournament entry #74
$ws.Cells.Item(75, 1).Value = 74
$ws.Range("B75").Value = 1271.848433724673
$ws.Range("C75").Value = "'"
$ws.Range("C75").Style = "Normal"
$ws.Range("D75").Value = "'"
$ws.Range("D75").Style = "Normal"
$ws.Range("E75").Value = "'"
$ws.Range("E75").Style = "Normal"
$ws.Range("F75").Value = 1285.225869219221

